$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    2.5303282475197002,
    2.54318725804539,
    2.55633262086624,
    2.5709707003347599,
    2.5851627489545401,
    2.59342596145909,
    2.60162425409395,
    2.6188078153578398,
    2.64240060943599,
    2.6688999520605501,
    2.6937108627031301,
    2.7058237660023599,
    2.7121811026100402,
    2.7157002466516298,
    2.71854095924345,
    2.7219876639975902,
    2.7239187288901001,
    2.72359693305981,
    2.7221034039346899,
    2.7183511365078799,
    2.7069797438571999,
    2.6947341612836602,
    2.6883715962836199,
    2.6848469037189302,
    2.6821331312997798,
    2.6815362260294902,
    2.68912185048121,
    2.6924721488391401,
    2.6920286680807801,
    2.68953240781914,
    2.6871961071242501,
    2.6848902085094499,
    2.68337590327658,
    2.68224870104561,
    2.6811214988146501,
    2.6799942965836898,
    2.67731031282108,
    2.6770548278194202,
    2.6793408213313299,
    2.68339216378439,
    2.68744350623744,
    2.6915580984716301,
    2.6954596007596399,
    2.6979860951591701,
    2.6987828441670101,
    2.7008907197673402,
    2.7036680410886,
    2.7068457489673499,
    2.71104098005046,
    2.71686448538347,
    2.72451014772663,
    2.7324519950528399,
    2.7407027219714499,
    2.74893798753603,
    2.75548833304325,
    2.7513246132226299,
    2.7456559443353599,
    2.73892923640391,
    2.7324451417298201,
    2.7281110445201202,
    2.7304319448593701,
    2.7393281660465498,
    2.75087738855221,
    2.76651371037135,
    2.7880652334660798,
    2.80939363193772,
    2.8403633086779201,
    2.8729774131315202,
    2.9039480715811501,
    2.9319539064101798,
    2.9601171756142901,
    2.9841894389635,
    3.0054419028049999,
    3.0259618093058802,
    3.04795993814056,
    3.06842853759477,
    3.0870683660279599,
    3.0962848190277898,
    3.0975469466399201,
    3.0969099137031,
    3.0950130483766598,
    3.0938670858438901,
    3.0917112326548302,
    3.0859578547313098,
    3.0735595811160898,
    3.0527246140500401,
    3.0266941751002099,
    2.9994306786688099,
    2.9693798857651701,
    2.9314636673908399,
    2.8855901033769702,
    2.8388720474887901,
    2.7919047972770499,
    2.7457355729117801,
    2.7019201856068999,
    2.6697865209429401,
    2.63696568675197,
    2.6039272488774601,
    2.5708888110029502,
    2.5378503731284399
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
